$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "team record" columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold/centered/bordered) from an existing header cell
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record values for every data row (2 through 43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 80   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 82   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
